# Trade #52 closed at 2026-02-17 13:29:04 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" rollups with the new trade
# totals, then appends the closed trade as a new row (53) to both the
# "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.47
$summary.Range("B4").Value = -2.52
$summary.Range("B5").Value = -0.97
$summary.Range("B6").Value = 52
$summary.Range("B7").Value = 20
$summary.Range("B9").Value = 38.46

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.47
$status.Range("D4").Value = 52
$status.Range("E4").Value = -2.52
$status.Range("F4").Value = -2.53
$status.Range("G4").Value = 38.46

# ---------------------------------------------------------------------
# New trade row (#52, spreadsheet row 53) appended to both the
# "All Trades" log and the per-strategy "MarketMaking" log.
# ---------------------------------------------------------------------
$newRow = 53

$tradeNumber   = 52
$tradeDate     = "2026-02-17"
$tradeTime     = "13:28:57"
$tradeStrategy = "MarketMaking"
$tradeSide     = "DOWN"
$entryPrice    = 0.11
$exitPrice     = 0.159574
$tradeStatus   = "CLOSED"
$pnlPct        = 45.0677
$pnlDollar     = 0.05
$capitalAfter  = 97.47
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$exitReason    = "early_exit"
$durationMin   = 0.12

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value  = $tradeNumber

    # The date column stores a plain "yyyy-mm-dd" text label (matching
    # every other row in the log), not a real Excel date serial. Force a
    # text number format before assigning so the literal string is kept
    # instead of being auto-converted to a date value, then drop back to
    # the default style so no stray formatting is left behind.
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = $tradeDate
    $ws.Cells.Item($newRow, 2).Style = "Normal"

    $ws.Cells.Item($newRow, 3).Value  = $tradeTime
    $ws.Cells.Item($newRow, 4).Value  = $tradeStrategy
    $ws.Cells.Item($newRow, 5).Value  = $tradeSide
    $ws.Cells.Item($newRow, 6).Value  = $entryPrice
    $ws.Cells.Item($newRow, 7).Value  = $exitPrice
    $ws.Cells.Item($newRow, 8).Value  = $tradeStatus
    $ws.Cells.Item($newRow, 9).Value  = $pnlPct
    $ws.Cells.Item($newRow, 10).Value = $pnlDollar
    $ws.Cells.Item($newRow, 11).Value = $capitalAfter
    $ws.Cells.Item($newRow, 12).Value = $entrySlippage
    $ws.Cells.Item($newRow, 13).Value = $exitSlippage
    $ws.Cells.Item($newRow, 14).Value = $confidence
    $ws.Cells.Item($newRow, 15).Value = $entryReason
    $ws.Cells.Item($newRow, 16).Value = $exitReason
    $ws.Cells.Item($newRow, 17).Value = $durationMin
}
